# Refresh cryptocurrency Price / Volume(1h) figures in the 'cryptos'
# worksheet, matching the GitHub Actions data pull on
# Sun Sep 29 20:37:49 UTC 2024. Also corrects the Stacks /
# FirstDigitalUSD row ordering (rows 38-39 swapped places).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: mark numeric-looking Price cells as Text so Excel does not
# coerce the assigned string into a floating point number (which would
# drop meaningful trailing zeros, e.g. '0.0620' -> 0.062).
$textProtectCells = @('D5', 'D6', 'D7', 'D8', 'D13', 'D18', 'D20', 'D25', 'D26', 'D28', 'D31', 'D34', 'D35', 'D36', 'D37', 'D38', 'D39', 'D40', 'D42', 'D43', 'D44', 'D45', 'D46', 'D49', 'D51')
foreach ($cellAddr in $textProtectCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

# Step 2: write the refreshed coin data.
$ws.Range('D2').Value = '65.926.93'
$ws.Range('E2').Value = '  +0.25%  '
$ws.Range('D3').Value = '2.663.69'
$ws.Range('E3').Value = '  -0.57%  '
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').Value = '599.55'
$ws.Range('E5').Value = '  -0.15%  '
$ws.Range('D6').Value = '159.90'
$ws.Range('E6').Value = '  +1.87%  '
$ws.Range('D7').Value = '0.646'
$ws.Range('E7').Value = '  +3.88%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('E9').Value = '  -2.87%  '
$ws.Range('E10').Value = '  +0.18%  '
$ws.Range('E11').Value = '  -0.23%  '
$ws.Range('E12').Value = '  +1.52%  '
$ws.Range('D13').Value = '29.11'
$ws.Range('E13').Value = '  -0.89%  '
$ws.Range('E14').Value = '  -1.25%  '
$ws.Range('D15').Value = '3.142.56'
$ws.Range('E15').Value = '  -0.53%  '
$ws.Range('D16').Value = '65.788.27'
$ws.Range('E16').Value = '  -0.94%  '
$ws.Range('D17').Value = '2.633.93'
$ws.Range('E17').Value = '  -1.67%  '
$ws.Range('D18').Value = '12.64'
$ws.Range('E18').Value = '  -2.22%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('D20').Value = '354.16'
$ws.Range('E20').Value = '  +0.41%  '
$ws.Range('E21').Value = '  -1.49%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E24').Value = '  +9.92%  '
$ws.Range('D25').Value = '0.0000114'
$ws.Range('E25').Value = '  +2.01%  '
$ws.Range('D26').Value = '9.68'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  +1.85%  '
$ws.Range('D28').Value = '579.24'
$ws.Range('E28').Value = '  +9.29%  '
$ws.Range('E29').Value = '  +1.57%  '
$ws.Range('E30').Value = '  -1.94%  '
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  -0.07%  '
$ws.Range('E32').Value = '  +0.84%  '
$ws.Range('D34').Value = '6.75'
$ws.Range('E34').Value = '  +4.21%  '
$ws.Range('D35').Value = '5.55'
$ws.Range('E35').Value = '  +0.75%  '
$ws.Range('D36').Value = '0.424'
$ws.Range('E36').Value = '  -0.12%  '
$ws.Range('D37').Value = '20.63'
$ws.Range('E37').Value = '  -0.61%  '
$ws.Range('B38').Value = 'Stacks'
$ws.Range('C38').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D38').Value = '1.97'
$ws.Range('E38').Value = '  +1.29%  '
$ws.Range('B39').Value = 'FirstDigitalUSD'
$ws.Range('C39').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D39').Value = '0.999'
$ws.Range('E39').Value = '  -0.10%  '
$ws.Range('D40').Value = '155.12'
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('E41').Value = '  +8.81%  '
$ws.Range('D42').Value = '161.70'
$ws.Range('E42').Value = '  -1.64%  '
$ws.Range('D43').Value = '4.12'
$ws.Range('E43').Value = '  -0.82%  '
$ws.Range('D44').Value = '0.0620'
$ws.Range('E44').Value = '  +1.18%  '
$ws.Range('D45').Value = '23.52'
$ws.Range('E45').Value = '  +2.68%  '
$ws.Range('D46').Value = '0.645'
$ws.Range('E46').Value = '  +0.31%  '
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('D49').Value = '19.76'
$ws.Range('E49').Value = '  -2.31%  '
$ws.Range('D50').Value = '0.0₆0246'
$ws.Range('E50').Value = '  -6.95%  '
$ws.Range('D51').Value = '0.820'
$ws.Range('E51').Value = '  +0.50%  '

# Step 3: restore the default cell style so the text-format tweak from
# step 1 doesn't leave a stray style behind on these cells.
foreach ($cellAddr in $textProtectCells) {
    $ws.Range($cellAddr).Style = "Normal"
}

